$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (e.g. "1.000" -> 1), so they stay as
# inline text strings, matching the source data feed format.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values (coin order, prices, and 1h volume %).
$ws.Range("D2").Value = "30.586.99"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.923.19"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "247.32"
$ws.Range("E5").Value = "  +2.73%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.4731"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "0.2916"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").Value = "0.06848"
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("D10").Value = "106.07"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").Value = "18.47"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("D12").Value = "1.930.13"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "0.07732"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "5.340"
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("D15").Value = "0.6733"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "289.36"
$ws.Range("E16").Value = "  -5.82%  "
$ws.Range("D17").Value = "30.611.82"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "0.000007650"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "12.98"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.569"
$ws.Range("E21").Value = "  +5.14%  "
$ws.Range("D22").Value = "2.180.81"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "6.497"
$ws.Range("D25").Value = "9.543"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").Value = "166.97"
$ws.Range("D27").Value = "20.79"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "2.134"
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("D29").Value = "0.1070"
$ws.Range("E29").Value = "  -3.32%  "
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("D31").Value = "4.206"
$ws.Range("E31").Value = "  +2.75%  "
$ws.Range("D32").Value = "4.075"
$ws.Range("E32").Value = "  +3.35%  "
$ws.Range("D33").Value = "0.05048"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "0.7348"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").Value = "1.148"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").Value = "0.02056"
$ws.Range("E36").Value = "  +4.72%  "
$ws.Range("D37").Value = "2.744"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "2.684"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "2.056"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "111.75"
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("D42").Value = "0.4469"
$ws.Range("E42").Value = "  +6.43%  "
$ws.Range("D43").Value = "0.8731"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D44").Value = "5.909"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "68.01"
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("D47").Value = "7.308"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "9.454"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").Value = "48.96"
$ws.Range("E49").Value = "  +14.22%  "
$ws.Range("E50").Value = "  +3.93%  "
$ws.Range("D51").Value = "35.30"
$ws.Range("E51").Value = "  +1.25%  "
